# Fills in the previously-empty (inlineStr placeholder) cells in columns D:H
# of Sheet1 with numeric values, per the "Remove empty space and add 0
# values" commit.
#
# Two groups of rows are affected:
#   1) Rows where D already has a budget cost but E/F/G/H were empty:
#      set E = -D (the time delta offsets the budget to zero) and
#      F = G = H = 0.
#   2) Rows where D/E were empty but F/G/H already had values:
#      set D = 0 and E = F (the time delta equals the actual cost,
#      since there is no budget).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsBudgetOnly = @(4, 6, 8, 10, 12, 13, 16, 24, 35, 43, 49, 51)
$rowsActualOnly = @(5, 7, 11, 17, 18, 19, 20, 21, 30, 31, 41, 44, 45, 46, 47, 48, 50, 52)

foreach ($r in $rowsBudgetOnly) {
    $budget = $ws.Range("D$r").Value2
    $ws.Range("E$r").Value2 = -$budget
    $ws.Range("F$r").Value2 = 0
    $ws.Range("G$r").Value2 = 0
    $ws.Range("H$r").Value2 = 0
}

foreach ($r in $rowsActualOnly) {
    $actual = $ws.Range("F$r").Value2
    $ws.Range("D$r").Value2 = 0
    $ws.Range("E$r").Value2 = $actual
}
